# Convert the "On Pilgrimage - March 1968" / "By Dorothy Day" heading block
# into a pandoc-style title block: a Title-styled paragraph with the title
# split word-by-word into runs, followed by an Authors-styled paragraph
# ("Dorothy Day") also split word-by-word into runs.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Paragraph 1 is the "On Pilgrimage - March 1968" heading (style Heading1),
# paragraph 2 is the "By Dorothy Day" byline (bold run, no named style).
$titlePara = $d.Paragraphs(1)
$authorPara = $d.Paragraphs(2)

# Sanity-check we are editing the paragraphs we expect before mutating.
# (Paragraph .Range.Text includes the trailing paragraph-mark character,
# so trim before comparing.)
if ($titlePara.Range.Text.TrimEnd() -ne "On Pilgrimage - March 1968") {
    throw "Unexpected paragraph 1 text: $($titlePara.Range.Text)"
}
if ($authorPara.Range.Text.TrimEnd() -ne "By Dorothy Day") {
    throw "Unexpected paragraph 2 text: $($authorPara.Range.Text)"
}

function Make-Run($text) {
    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    return "<w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r>"
}

$titleWords = @("On", " ", "Pilgrimage", " ", "-", " ", "March", " ", "1968")
$authorWords = @("Dorothy", " ", "Day")

$titleRuns = ($titleWords | ForEach-Object { Make-Run $_ }) -join ""
$authorRuns = ($authorWords | ForEach-Object { Make-Run $_ }) -join ""

$newXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Title`"/></w:pPr>$titleRuns</w:p>" +
          "<w:p $wNs><w:pPr><w:pStyle w:val=`"Authors`"/></w:pPr>$authorRuns</w:p>"

# Replace both paragraphs (title heading + byline) in one shot so the old
# runs/paragraph styling are fully discarded and replaced by the new
# pandoc-style title block.
$range = $d.Range($titlePara.Range.Start, $authorPara.Range.End)
$range.InsertXML($newXml) | Out-Null
